$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2864.3333
$ws.Range("I62").Value = 2932.2222
$ws.Range("K62").Value = 2932.2222
$ws.Range("M62").Value = -2308.2222
# Row 65
$ws.Range("H65").Value = 2864.3333
$ws.Range("I65").Value = 2932.2222
$ws.Range("K65").Value = 14661.111
$ws.Range("M65").Value = -11541.111
# Row 87
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 111
$ws.Range("H111").Value = 1749.8667
$ws.Range("J111").Value = 1504.8334
$ws.Range("L111").Value = 4514.5002
$ws.Range("N111").Value = -10648.5002
# Row 116
$ws.Range("H116").Value = 4142.857
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4142.857
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").Value = 4142.857
$ws.Range("N116").Value = -11026.857
# Row 125
$ws.Range("H125").Value = 1466.3334
$ws.Range("J125").Value = 1724.5
$ws.Range("L125").Value = 15520.5
$ws.Range("N125").Value = -20440.5
# Row 137
$ws.Range("H137").Value = 24393062
$ws.Range("I137").Value = 1202.72
$ws.Range("J137").Value = 62505344
$ws.Range("K137").Value = 3608.16
$ws.Range("L137").Value = 187516032
$ws.Range("M137").Value = -1058.16
$ws.Range("N137").Value = -187521132
# Row 138
$ws.Range("H138").Value = 2946.9092
$ws.Range("I138").Value = 1211.5172
$ws.Range("J138").Value = 3665.8572
$ws.Range("K138").Value = 3634.5516
$ws.Range("L138").Value = 10997.5716
$ws.Range("M138").Value = 1505.4484
$ws.Range("N138").Value = -21277.5716
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 38410.27
$ws.Range("I32").Value = 37509.51
$ws.Range("J32").Value = 46742.25
$ws.Range("K32").Value = 37509.51
$ws.Range("L32").Value = 46742.25
$ws.Range("M32").Value = -37222.51
$ws.Range("N32").Value = -47316.25
# Row 37
$ws.Range("H37").Value = 9790.444
$ws.Range("I37").Value = 2086.8
$ws.Range("J37").Value = 19420
$ws.Range("K37").Value = 2086.8
$ws.Range("L37").Value = 19420
$ws.Range("M37").Value = -1813.8
$ws.Range("N37").Value = -19966
# Row 61
$ws.Range("H61").Value = 1926.8462
$ws.Range("I61").Value = 1786.8695
$ws.Range("K61").Value = 1786.8695
$ws.Range("M61").Value = -1574.8695
# Row 102
$ws.Range("H102").Value = 2603.3333
$ws.Range("I102").Value = 2603.3333
$ws.Range("K102").Value = 2603.3333
$ws.Range("M102").Value = -981.3332999999998
# Row 110
$ws.Range("H110").Value = 1384.2941
$ws.Range("I110").Value = 1621
$ws.Range("J110").Value = 615
$ws.Range("K110").Value = 1621
$ws.Range("L110").Value = 615
$ws.Range("M110").Value = 424
$ws.Range("N110").Value = -4705
# Row 136
$ws.Range("H136").Value = 1926.8462
$ws.Range("I136").Value = 1786.8695
$ws.Range("K136").Value = 5360.6085
$ws.Range("M136").Value = -2810.6085

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 31
$ws.Range("H31").Value = 12000
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 86
$ws.Range("H86").Value = 2522.762
$ws.Range("J86").Value = 1989
$ws.Range("L86").Value = 1989
$ws.Range("N86").Value = -4235
# Row 89
$ws.Range("H89").Value = 2522.762
$ws.Range("J89").Value = 1989
$ws.Range("L89").Value = 9945
$ws.Range("N89").Value = -21177
# Row 94
$ws.Range("H94").Value = 386.15384
$ws.Range("I94").Value = 385.9091
$ws.Range("J94").Value = 387.5
$ws.Range("K94").Value = 385.9091
$ws.Range("L94").Value = 387.5
$ws.Range("M94").Value = 65.09089999999998
$ws.Range("N94").Value = -1289.5
# Row 104
$ws.Range("H104").Value = 33552.6
$ws.Range("J104").Value = 33552.6
$ws.Range("L104").Value = 33552.6
$ws.Range("N104").Value = -40540.6
# Row 133
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 49000
$ws.Range("L133").Value = 49000
$ws.Range("N133").Value = -59120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 27786152
$ws.Range("I31").Value = 9028.764999999999
$ws.Range("J31").Value = 52639370
$ws.Range("K31").Value = 9028.764999999999
$ws.Range("L31").Value = 52639370
$ws.Range("M31").Value = -8733.764999999999
$ws.Range("N31").Value = -52639960
# Row 34
$ws.Range("H34").Value = 27786152
$ws.Range("I34").Value = 9028.764999999999
$ws.Range("J34").Value = 52639370
$ws.Range("K34").Value = 9028.764999999999
$ws.Range("L34").Value = 52639370
$ws.Range("M34").Value = -8826.764999999999
$ws.Range("N34").Value = -52639774
# Row 41
$ws.Range("H41").Value = 4300
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
# Row 58
$ws.Range("H58").Value = 2054.3
$ws.Range("I58").Value = 1900.6666
$ws.Range("J58").Value = 2668.8333
$ws.Range("K58").Value = 1900.6666
$ws.Range("L58").Value = 2668.8333
$ws.Range("M58").Value = -1697.6666
$ws.Range("N58").Value = -3074.8333
# Row 94
$ws.Range("H94").Value = 142858610
$ws.Range("I94").Value = 1000000000
$ws.Range("J94").Value = 1704.3334
$ws.Range("K94").Value = 1000000000
$ws.Range("L94").Value = 1704.3334
$ws.Range("M94").Value = -999999549
$ws.Range("N94").Value = -2606.3334
# Row 132
$ws.Range("H132").Value = 3908078.5
$ws.Range("I132").Value = 1572
$ws.Range("K132").Value = 4716
$ws.Range("M132").Value = -2186
# Row 134
$ws.Range("H134").Value = 5834.773
$ws.Range("I134").Value = 7279.0625
$ws.Range("J134").Value = 1983.3334
$ws.Range("K134").Value = 21837.1875
$ws.Range("L134").Value = 5950.0002
$ws.Range("M134").Value = -19302.1875
$ws.Range("N134").Value = -11020.0002
# Row 136
$ws.Range("H136").Value = 2054.3
$ws.Range("I136").Value = 1900.6666
$ws.Range("J136").Value = 2668.8333
$ws.Range("K136").Value = 5701.9998
$ws.Range("L136").Value = 8006.499899999999
$ws.Range("M136").Value = -3151.9998
$ws.Range("N136").Value = -13106.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 3866.6667
$ws.Range("J39").Value = 3866.6667
$ws.Range("L39").Value = 11600.0001
$ws.Range("N39").Value = -12188.0001
# Row 55
$ws.Range("H55").Value = 2225.3333
$ws.Range("J55").Value = 2225.3333
$ws.Range("L55").Value = 6675.999899999999
$ws.Range("N55").Value = -7029.999899999999
# Row 107
$ws.Range("H107").Value = 522.2632
$ws.Range("I107").Value = 586.6667
$ws.Range("J107").Value = 492.53845
$ws.Range("K107").Value = 1760.0001
$ws.Range("L107").Value = 1477.61535
$ws.Range("M107").Value = 159.9999
$ws.Range("N107").Value = -5317.61535
# Row 131
$ws.Range("H131").Value = 720.37
$ws.Range("I131").Value = 292
$ws.Range("J131").Value = 795.9647
$ws.Range("K131").Value = 876
$ws.Range("L131").Value = 2387.8941
$ws.Range("M131").Value = 4164
$ws.Range("N131").Value = -12467.8941

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 372.14285
$ws.Range("J107").Value = 159.66667
$ws.Range("L107").Value = 159.66667
$ws.Range("N107").Value = -3999.66667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 535.9
$ws.Range("I22").Value = 727.8
$ws.Range("J22").Value = 344
$ws.Range("K22").Value = 727.8
$ws.Range("L22").Value = 344
$ws.Range("M22").Value = -432.8
$ws.Range("N22").Value = -934
# Row 27
$ws.Range("H27").Value = 535.9
$ws.Range("I27").Value = 727.8
$ws.Range("J27").Value = 344
$ws.Range("K27").Value = 727.8
$ws.Range("L27").Value = 344
$ws.Range("M27").Value = -620.8
$ws.Range("N27").Value = -558
# Row 46
$ws.Range("H46").Value = 1500.7894
$ws.Range("I46").Value = 875
$ws.Range("J46").Value = 1574.4117
$ws.Range("K46").Value = 875
$ws.Range("L46").Value = 1574.4117
$ws.Range("M46").Value = -687
$ws.Range("N46").Value = -1950.4117
# Row 55
$ws.Range("H55").Value = 292.2619
$ws.Range("I55").Value = 313.33334
$ws.Range("J55").Value = 271.1905
$ws.Range("K55").Value = 313.33334
$ws.Range("L55").Value = 271.1905
$ws.Range("M55").Value = -140.33334
$ws.Range("N55").Value = -617.1904999999999
# Row 122
$ws.Range("H122").Value = 5117.3335
$ws.Range("I122").Value = 9402
$ws.Range("K122").Value = 28206
$ws.Range("M122").Value = -25756

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 338.46667
$ws.Range("I100").Value = 256.5
$ws.Range("K100").Value = 513
$ws.Range("M100").Value = 28
# Row 136
$ws.Range("H136").Value = 2230.0356
$ws.Range("I136").Value = 2226.04
$ws.Range("J136").Value = 2263.3333
$ws.Range("K136").Value = 6678.12
$ws.Range("L136").Value = 6789.999899999999
$ws.Range("M136").Value = -4128.12
$ws.Range("N136").Value = -11889.9999
